# Apply updated "dSF" (column F) values to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 2
    4  = -1
    5  = 1
    6  = -3
    7  = 4
    9  = -3
    10 = 1
    11 = 2
    12 = -2
    13 = -5
    14 = 4
    15 = 10
    16 = -3
    18 = -2
    20 = 3
    21 = -2
    22 = -1
    23 = -2
    25 = -1
    26 = -1
    27 = -1
    28 = -1
    29 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
